$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(10.44645118713379, -1.072464942932129, -7.741544246673584),
  @(-12.54131031036377, -6.995788097381592, 4.960752010345459),
  @(1.563988208770752, -0.6327500343322754, 7.09266471862793),
  @(1.386142730712891, -3.08348274230957, 6.775681018829346),
  @(1.111974596977234, -5.867249965667725, 13.13678932189941),
  @(-4.935175895690918, 8.790708541870117, 6.229645729064941),
  @(16.58664321899414, -65.17784118652344, 6.808780193328857),
  @(-8.056964874267578, 23.69245338439941, -1.875571012496948),
  @(-4.308743476867676, 2.427361011505127, 6.903748512268066),
  @(-2.216989278793335, -25.4277229309082, 2.421021461486816),
  @(-0.3425538539886474, -11.7038745880127, 31.56782150268555),
  @(-1.816783547401428, 8.313332557678223, 9.763429641723633),
  @(2.05627965927124, -54.87781143188477, 29.65373039245605),
  @(-26.46129989624023, 36.53782653808594, 2.705925941467285),
  @(-6.459963798522949, 9.631237030029297, -4.527087211608887),
  @(-2.323664665222168, -6.97331714630127, 3.054933786392212),
  @(23.92743682861328, 6.633898258209229, 20.45568084716797),
  @(12.88255214691162, 13.76539325714111, 5.36094856262207),
  @(-39.15726470947266, -50.35159301757812, 59.05854797363281),
  @(-29.70075225830078, 18.98210144042969, -6.946440696716309),
  @(-3.437598705291748, 8.676052093505859, -6.736623287200928),
  @(36.55035400390625, -4.483262062072754, -3.189533472061157),
  @(25.18490791320801, 10.66421031951904, 36.16248321533203),
  @(11.78367233276367, 19.32002067565918, 14.97337532043457),
  @(-10.2064151763916, -54.4849967956543, 45.01205825805664),
  @(-5.304520606994629, 4.910325050354004, -39.37523651123047),
  @(10.93332672119141, 14.53017807006836, -3.883467674255371),
  @(23.18131637573243, -33.343994140625, -4.427485942840576),
  @(17.88149261474609, -17.56607437133789, -2.040470600128174),
  @(-4.381585597991943, 13.65173721313477, 6.552346229553223)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
